# Apply updates described by the diff to sheet1 of the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update footnote texts in row 8 (A8, B8, C8) ---
# Order of assignment controls the shared-string table order produced by the
# engine; match the target string order: B8, then C8, then A8.
$ws.Range("B8").Value = "*по данным МЦР КР"
$ws.Range("C8").Value = "*according to the MDD KR"
$ws.Range("A8").Value = "*КР СӨМ маалыматтары  боюнча"

# --- Add new column O with 2023 data ---
# Copy column N (rows 3-7) into column O to inherit matching styles, then
# overwrite the cell values that differ for the new year.
$ws.Range("N3:N7").Copy()
$ws.Range("O3:O7").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("O4").Value = 2023
$ws.Range("O5").Value = 99
$ws.Range("O6").Value = 98.9
$ws.Range("O7").Value = 98.8

# --- Adjust column widths for columns A:C ---
# The engine's ColumnWidth setter adds a constant ~0.8333 offset relative to
# the width value actually persisted in the XML, so back that off here to
# land on a stored width of exactly 38.
$ws.Range("A1:C1").ColumnWidth = 37.166666666666664

# --- Reset the active selection back to A1 ---
$ws.Range("A1").Select()
